$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "AL QUR'AN EDISI TAHLILAN 30 Juz + Doa Tahlil | Pengganti Buku Yasin | Al Aqeel A6 Pastel HVS Edisi Tahlilan"
$ws.Range("B22").Value = "Custom sisipan 1 hal"
$ws.Range("C22").Value = 1250

$ws.Range("B22").Borders.Item(7).LineStyle = 1
$ws.Range("B22").Borders.Item(7).Weight = 2
$ws.Range("B22").Borders.Item(10).LineStyle = 1
$ws.Range("B22").Borders.Item(10).Weight = 2
